# Simplify the attendance-template header sheet:
#  - Rename the first header cell from the Arabic name prompt to "*Your Name"
#  - Clear out the old "answers" label in B1 (keep its style/formatting)
#  - Drop the whole "select" / default-value column (C) entirely
#  - Replace the "left empty in answers" note with randomization guidance
#  - Re-flow column A to take up the freed width, and move the selection
#    down to A2 (matches random_select enhancement for quizzes)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new English header text; B1 becomes a blank (but still styled) cell
$ws.Range("A1").Value = "*Your Name"
$ws.Range("B1").ClearContents()

# Column C ("select" header, the literal 0 default and anything else in it)
# is no longer needed now that there's no dedicated answers column.
$ws.Columns("C").Delete()

# Row 3: clarify how #TEXT#/#NUMBER# should be chosen
$ws.Range("A3").Value = "Choose types whether #TEXT# or #NUMBER#"

# Column A widens to take over the space freed by the removed column
$ws.Columns("A").ColumnWidth = 40.6667

# Move the active selection to A2
$ws.Range("A2").Select()
